# Applies odds/value updates to Sheet1 of the FlashScore weekly games workbook
# as described by the commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.62
$ws.Range("G3").Value = 2.4
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 3.1
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.57
$ws.Range("AG3").Value = 351
$ws.Range("AO3").Value = 13
$ws.Range("AR3").Value = 81
$ws.Range("AX3").Value = 19
$ws.Range("AZ3").Value = 67
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 4.75
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.57
$ws.Range("BD4").Value = 151
$ws.Range("G8").Value = 4.1
$ws.Range("I8").Value = 1.75
$ws.Range("J8").Value = 4.5
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.93
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("W8").Value = 12
$ws.Range("X8").Value = 21
$ws.Range("Z8").Value = 41
$ws.Range("AT8").Value = 2.75
$ws.Range("AY8").Value = 21
$ws.Range("AZ8").Value = 34
$ws.Range("BB8").Value = 151
$ws.Range("H11").Value = 4.2
$ws.Range("O12").Value = 1.18
$ws.Range("P12").Value = 4.5
$ws.Range("Q12").Value = 1.65
$ws.Range("R12").Value = 2.2
$ws.Range("G13").Value = 4.2
$ws.Range("H13").Value = 4.33
$ws.Range("I13").Value = 1.67
$ws.Range("J13").Value = 4.33
$ws.Range("K13").Value = 2.6
$ws.Range("L13").Value = 2.2
$ws.Range("U13").Value = 1.44
$ws.Range("V13").Value = 2.63
$ws.Range("W13").Value = 19
$ws.Range("X13").Value = 26
$ws.Range("Z13").Value = 41
$ws.Range("AB13").Value = 26
$ws.Range("AC13").Value = 23
$ws.Range("AD13").Value = 9
$ws.Range("AJ13").Value = 9
$ws.Range("AN13").Value = 6.5
$ws.Range("AP13").Value = 21
$ws.Range("AQ13").Value = 51
$ws.Range("AW13").Value = 4.33
$ws.Range("AX13").Value = 8.5
$ws.Range("AY13").Value = 15
$ws.Range("AZ13").Value = 23
$ws.Range("M23").Value = 1.05
$ws.Range("O23").Value = 1.41
$ws.Range("P23").Value = 2.62
$ws.Range("V23").Value = 1.69
$ws.Range("M24").Value = 1.05
$ws.Range("O24").Value = 1.37
$ws.Range("U24").Value = 1.87
$ws.Range("V24").Value = 1.77
$ws.Range("O28").Value = 1.25
$ws.Range("P28").Value = 3.75
$ws.Range("Q28").Value = 1.93
$ws.Range("R28").Value = 1.93
$ws.Range("G30").Value = 2.1
$ws.Range("I30").Value = 3.6
$ws.Range("L30").Value = 4
$ws.Range("M30").Value = 1.05
$ws.Range("N30").Value = 9
$ws.Range("O30").Value = 1.27
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("W30").Value = 7.5
$ws.Range("Z30").Value = 19
$ws.Range("AL30").Value = 29
$ws.Range("AP30").Value = 21
$ws.Range("BA30").Value = 81
